$d = $word.ActiveDocument

$replacements = @(
    @("2023-11-24 Friday", "2023-11-25 Saturday"),
    @("16×77=", "70×12="),
    @("60×49=", "72×55="),
    @("71×96=", "14×62="),
    @("18×12=", "48×21="),
    @("86×85=", "46×15="),
    @("49×55=", "87×94="),
    @("79×44=", "29×15="),
    @("94×39=", "85×77="),
    @("83×15=", "87×69="),
    @("79×96=", "74×57="),
    @("19×97=", "16×59="),
    @("56×88=", "13×26="),
    @("29×88=", "64×96="),
    @("73×88=", "47×18="),
    @("94×26=", "11×49="),
    @("50×69=", "84×90="),
    @("58×44=", "62×84="),
    @("78×26=", "12×42="),
    @("94×50=", "21×78="),
    @("41×89=", "24×50="),
    @("32×49=", "50×61="),
    @("48×74=", "36×84="),
    @("18×57=", "94×95="),
    @("82×75=", "16×30="),
    @("84×97=", "74×58=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
